$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.129.55"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.978.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.14"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.08"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.518"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.978.84"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.145"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.02"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000226"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.16"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.470.95"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.90"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.178.08"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.983.96"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "447.80"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.15"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.680"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.30"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.98"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.17"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.39"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.87"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.68"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.06%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.11"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.04"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.108"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0816"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.48%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.75"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.05"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.04"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.122"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +9.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.83"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "393.61"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.28"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0349"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.266"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.685.21"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.08"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.13"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.20%  "
